# Releases y escenarios para planeamiento
# Adds SPRINT (column J) values to the visible backlog items on the
# FEATURES sheet and filters the list down to sprint SP1 via AutoFilter.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("FEATURES")

# Assign sprint values to the rows that currently belong to EN001
# (i.e. the rows that are visible under the existing EN001 filter).
$ws.Range("J4").Value  = "SP1"
$ws.Range("J6").Value  = "SP2"
$ws.Range("J7").Value  = "SP2"
$ws.Range("J31").Value = "SP2"
$ws.Range("J33").Value = "SP2"
$ws.Range("J38").Value = "SP2"
$ws.Range("J39").Value = "SP2"
$ws.Range("J41").Value = "SP2"
$ws.Range("J42").Value = "SP1"
$ws.Range("J43").Value = "SP1"
$ws.Range("J44").Value = "SP2"
$ws.Range("J45").Value = "SP1"
$ws.Range("J46").Value = "SP2"
$ws.Range("J66").Value = "SP2"
$ws.Range("J75").Value = "SP1"
$ws.Range("J76").Value = "SP1"
$ws.Range("J77").Value = "SP1"

# Apply an AutoFilter on the SPRINT column (J, field 10) restricted to
# SP1, which hides every row whose sprint is SP2 (or blank).
[void]$ws.Range("A3:J90").AutoFilter(10, @("SP1"), 7)

# Update the current selection to match the author's final position.
[void]$ws.Range("E75").Select()
